$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (all Price/Volume cells are stored as plain text
# in the source data, e.g. "406.50" or "  -1.29%  ", so force Text format before
# writing to stop Excel from re-interpreting numeric-looking strings as numbers
# (which would drop formatting like trailing zeros).
$updates = [ordered]@{
  'D2' = '61.446.91'
  'E2' = '  -1.29%  '
  'D3' = '3.374.99'
  'E3' = '  -2.35%  '
  'E4' = '  +0.07%  '
  'D5' = '406.50'
  'E5' = '  -1.88%  '
  'D6' = '134.16'
  'E6' = '  +8.13%  '
  'E7' = '  +0.05%  '
  'E8' = '  +0.07%  '
  'D9' = '0.669'
  'E9' = '  +1.33%  '
  'E10' = '  -6.59%  '
  'D11' = '42.63'
  'E11' = '  +3.04%  '
  'E12' = '  -1.04%  '
  'D13' = '3.895.63'
  'E13' = '  -2.78%  '
  'D14' = '8.41'
  'E14' = '  -1.24%  '
  'D15' = '19.69'
  'E15' = '  -0.36%  '
  'D16' = '3.382.37'
  'E16' = '  -2.07%  '
  'D17' = '61.456.61'
  'E17' = '  -1.14%  '
  'E18' = '  -1.58%  '
  'D19' = '10.99'
  'E19' = '  -2.52%  '
  'E20' = '  -4.32%  '
  'E21' = '  -4.56%  '
  'D22' = '85.05'
  'E22' = '  +3.41%  '
  'D23' = '314.87'
  'E23' = '  +0.08%  '
  'D24' = '12.81'
  'E24' = '  -1.88%  '
  'E25' = '  -1.72%  '
  'D26' = '4.79'
  'E26' = '  +11.54%  '
  'D27' = '8.33'
  'E27' = '  +5.43%  '
  'D28' = '29.51'
  'E28' = '  -5.31%  '
  'D29' = '7.62'
  'E29' = '  -2.62%  '
  'E30' = '  +0.56%  '
  'E31' = '  -1.43%  '
  'D32' = '2.62'
  'E32' = '  +1.60%  '
  'D33' = '11.35'
  'E33' = '  -1.95%  '
  'E34' = '  -0.13%  '
  'D35' = '40.75'
  'E35' = '  -4.11%  '
  'E36' = '  -0.54%  '
  'D37' = '51.87'
  'E37' = '  -0.92%  '
  'D38' = '0.999'
  'E38' = '  +0.12%  '
  'E39' = '  -2.51%  '
  'D40' = '2.93'
  'E40' = '  -2.47%  '
  'D41' = '139.20'
  'E42' = '  -2.58%  '
  'E43' = '  -1.04%  '
  'D44' = '0.295'
  'E44' = '  +3.14%  '
  'D45' = '4.03'
  'E45' = '  +3.48%  '
  'D46' = '16.72'
  'E46' = '  -3.06%  '
  'E47' = '  +0.53%  '
  'D48' = '21.36'
  'E48' = '  -4.44%  '
  'D49' = '2.121.01'
  'E49' = '  -4.02%  '
  'D50' = '2.29'
  'E50' = '  -5.15%  '
  'E51' = '  +1.00%  '
}

foreach ($cellRef in $updates.Keys) {
  $cell = $ws.Range($cellRef)
  $cell.NumberFormat = "@"
  $cell.Value = $updates[$cellRef]
}
